$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column G, matching the style of column F's header
$ws.Range("G1").Value = "E"
$ws.Range("G1").Style = $ws.Range("F1").Style

# Update data rows 2-10: columns B-F get a "_<n>" suffix (where n = row-1),
# and new column G gets "e<n>"
$ws.Range("B2").Value = "M_1"
$ws.Range("C2").Value = "A_1"
$ws.Range("D2").Value = "N_1"
$ws.Range("E2").Value = "D1_1"
$ws.Range("F2").Value = "D2_1"
$ws.Range("G2").Value = "e1"

$ws.Range("B3").Value = "M_2"
$ws.Range("C3").Value = "A_2"
$ws.Range("D3").Value = "N_2"
$ws.Range("E3").Value = "D1_2"
$ws.Range("F3").Value = "D2_2"
$ws.Range("G3").Value = "e2"

$ws.Range("B4").Value = "M_3"
$ws.Range("C4").Value = "A_3"
$ws.Range("D4").Value = "N_3"
$ws.Range("E4").Value = "D1_3"
$ws.Range("F4").Value = "D2_3"
$ws.Range("G4").Value = "e3"

$ws.Range("B5").Value = "M_4"
$ws.Range("C5").Value = "A_4"
$ws.Range("D5").Value = "N_4"
$ws.Range("E5").Value = "D1_4"
$ws.Range("F5").Value = "D2_4"
$ws.Range("G5").Value = "e4"

$ws.Range("B6").Value = "M_5"
$ws.Range("C6").Value = "A_5"
$ws.Range("D6").Value = "N_5"
$ws.Range("E6").Value = "D1_5"
$ws.Range("F6").Value = "D2_5"
$ws.Range("G6").Value = "e5"

$ws.Range("B7").Value = "M_6"
$ws.Range("C7").Value = "A_6"
$ws.Range("D7").Value = "N_6"
$ws.Range("E7").Value = "D1_6"
$ws.Range("F7").Value = "D2_6"
$ws.Range("G7").Value = "e6"

$ws.Range("B8").Value = "M_7"
$ws.Range("C8").Value = "A_7"
$ws.Range("D8").Value = "N_7"
$ws.Range("E8").Value = "D1_7"
$ws.Range("F8").Value = "D2_7"
$ws.Range("G8").Value = "e7"

$ws.Range("B9").Value = "M_8"
$ws.Range("C9").Value = "A_8"
$ws.Range("D9").Value = "N_8"
$ws.Range("E9").Value = "D1_8"
$ws.Range("F9").Value = "D2_8"
$ws.Range("G9").Value = "e8"

$ws.Range("B10").Value = "M_9"
$ws.Range("C10").Value = "A_9"
$ws.Range("D10").Value = "N_9"
$ws.Range("E10").Value = "D1_9"
$ws.Range("F10").Value = "D2_9"
$ws.Range("G10").Value = "e9"

# Apply the same style used in columns B-F to the new column G cells for rows 2-10
$ws.Range("G2:G10").Style = $ws.Range("F2:F10").Style

# Update the active cell selection to D15
$ws.Range("D15").Select()
